$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.261.85"
Set-TextValue "E2" "  +1.08%  "
Set-TextValue "D3" "1.919.70"
Set-TextValue "E3" "  +0.67%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "E5" "  +1.42%  "
Set-TextValue "D6" "244.31"
Set-TextValue "E6" "  +1.11%  "
Set-TextValue "E7" "  +0.12%  "
Set-TextValue "D8" "0.3248"
Set-TextValue "E8" "  +3.03%  "
Set-TextValue "D9" "26.97"
Set-TextValue "E9" "  +2.68%  "
Set-TextValue "D10" "0.07246"
Set-TextValue "E10" "  +4.93%  "
Set-TextValue "D11" "0.7904"
Set-TextValue "E11" "  +7.25%  "
Set-TextValue "D12" "0.08091"
Set-TextValue "E12" "  +1.32%  "
Set-TextValue "D13" "1.906.27"
Set-TextValue "E13" "  -0.08%  "
Set-TextValue "D14" "5.415"
Set-TextValue "E14" "  +4.40%  "
Set-TextValue "D15" "93.95"
Set-TextValue "E15" "  +1.17%  "
Set-TextValue "D16" "30.267.85"
Set-TextValue "E16" "  +1.16%  "
Set-TextValue "D17" "14.22"
Set-TextValue "E17" "  +2.00%  "
Set-TextValue "D18" "6.060"
Set-TextValue "E18" "  +3.48%  "
Set-TextValue "D19" "250.13"
Set-TextValue "E19" "  +1.97%  "
Set-TextValue "D20" "0.000007843"
Set-TextValue "E20" "  +1.45%  "
Set-TextValue "D21" "2.169.38"
Set-TextValue "E21" "  +0.82%  "
Set-TextValue "D22" "8.229"
Set-TextValue "E22" "  +20.84%  "
Set-TextValue "E23" "  +0.12%  "
Set-TextValue "D24" "1.001"
Set-TextValue "E24" "  +0.11%  "
Set-TextValue "D25" "0.1669"
Set-TextValue "E25" "  +18.32%  "
Set-TextValue "D26" "9.490"
Set-TextValue "E26" "  +3.31%  "
Set-TextValue "D27" "167.79"
Set-TextValue "E27" "  +0.09%  "
Set-TextValue "D28" "19.02"
Set-TextValue "E28" "  +0.73%  "
Set-TextValue "D29" "2.161"
Set-TextValue "E29" "  +6.58%  "
Set-TextValue "D30" "1.388"
Set-TextValue "E30" "  +1.97%  "
Set-TextValue "E31" "  +2.58%  "
Set-TextValue "D32" "4.338"
Set-TextValue "E32" "  +0.93%  "
Set-TextValue "D33" "0.05754"
Set-TextValue "E33" "  +5.48%  "
Set-TextValue "D34" "4.150"
Set-TextValue "E34" "  +1.68%  "
Set-TextValue "D35" "1.297"
Set-TextValue "E35" "  +2.91%  "
Set-TextValue "D36" "0.7504"
Set-TextValue "E36" "  +3.11%  "
Set-TextValue "B37" "Frax"
Set-TextValue "C37" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D37" "1.001"
Set-TextValue "E37" "  +0.18%  "
Set-TextValue "B38" "HuobiToken"
Set-TextValue "C38" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D38" "2.731"
Set-TextValue "E38" "  +0.45%  "
Set-TextValue "D39" "0.01961"
Set-TextValue "E39" "  +2.11%  "
Set-TextValue "D40" "2.823"
Set-TextValue "E40" "  +1.50%  "
Set-TextValue "D41" "0.4566"
Set-TextValue "E41" "  +3.46%  "
Set-TextValue "D42" "74.28"
Set-TextValue "E42" "  +2.91%  "
Set-TextValue "D43" "5.978"
Set-TextValue "E43" "  -2.61%  "
Set-TextValue "D44" "0.8518"
Set-TextValue "E44" "  +2.10%  "
Set-TextValue "D45" "1.930"
Set-TextValue "E45" "  +3.34%  "
Set-TextValue "E46" "  +0.15%  "
Set-TextValue "D47" "103.60"
Set-TextValue "E47" "  +3.20%  "
Set-TextValue "B48" "Maker"
Set-TextValue "C48" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D48" "1.032.25"
Set-TextValue "E48" "  +4.60%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "9.983"
Set-TextValue "E49" "  +2.85%  "
Set-TextValue "B50" "SynthetixNetwork"
Set-TextValue "C50" "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue "D50" "3.110"
Set-TextValue "E50" "  +12.62%  "
Set-TextValue "D51" "7.635"
Set-TextValue "E51" "  +1.55%  "
